$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 53. This shifts:
#   old A53 "US Census Bureau, 2007"  -> A54
#   old A54 url (hyperlink)           -> A55
#   old A55 "" (empty)                -> A56
#   old A58 "USCB"                    -> A59
#   old A59 (long citation)           -> A60
$ws.Rows(53).Insert()

# The URL text (now at A55, still styled/linked as a hyperlink) needs to move
# down to A56 as plain text, while A55 becomes the blank cell.
$urlText = $ws.Range("A55").Value2

# Remove the (now stale) hyperlink object entirely - the final layout has no
# hyperlink on the sheet at all.
if ($ws.Range("A55").Hyperlinks.Count -gt 0) {
    $ws.Range("A55").Hyperlinks.Delete()
}
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# Clear the old hyperlink-styled cell so it becomes the blank row, matching
# the plain "source" styling used by the other blank lines in this block.
$ws.Range("A55").ClearContents()
$ws.Range("A55").Font.Underline = $false
$ws.Range("A55").Font.Italic = $true
$ws.Range("A55").Font.ColorIndex = 0

# Put the URL text at A56 as plain text, formatted with the "source" style
# (italic, same as the surrounding source/footnote lines) instead of the
# hyperlink style.
$ws.Range("A56").Value = $urlText
$ws.Range("A56").Font.Italic = $true
$ws.Range("A56").Font.Underline = $false
$ws.Range("A56").Font.ColorIndex = 0

# Row 60 (previously the long citation, shifted down from the old A59) is
# replaced with a second "USCB" line, styled like the other source lines.
$ws.Range("A60").Value = "USCB"
$ws.Range("A60").Font.Italic = $true
$ws.Range("A60").Font.Bold = $false
$ws.Range("A60").Font.Underline = $false
